$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "01:00:00" values from columns K and N for rows 5 through 15
for ($r = 5; $r -le 15; $r++) {
    $ws.Cells.Item($r, 11).Value = ""
    $ws.Cells.Item($r, 14).Value = ""
}
